$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$ws.Range("B2").Value = 40
$ws.Range("D2").Value = 0

$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 0
